$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "28.348.43"
Set-TextCell "E2" "  +1.46%  "

Set-TextCell "D3" "1.823.65"
Set-TextCell "E3" "  +2.56%  "

Set-TextCell "D5" "317.24"
Set-TextCell "E5" "  +0.47%  "

Set-TextCell "D7" "0.5331"
Set-TextCell "E7" "  -0.98%  "

Set-TextCell "D8" "0.4043"
Set-TextCell "E8" "  +7.43%  "

Set-TextCell "D9" "0.07594"
Set-TextCell "E9" "  +2.15%  "

Set-TextCell "D10" "41.84"
Set-TextCell "E10" "  +0.54%  "

Set-TextCell "E11" "  +1.30%  "

Set-TextCell "D12" "6.315"
Set-TextCell "E12" "  +4.09%  "

Set-TextCell "E13" "  +0.02%  "

Set-TextCell "D14" "7.610"
Set-TextCell "E14" "  +5.82%  "

Set-TextCell "D15" "20.82"
Set-TextCell "E15" "  +2.03%  "

Set-TextCell "D16" "1.829.66"
Set-TextCell "E16" "  +2.91%  "

Set-TextCell "B17" "Litecoin"
Set-TextCell "C17" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D17" "89.37"
Set-TextCell "E17" "  +1.53%  "

Set-TextCell "B18" "ShibaInu"
Set-TextCell "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D18" "0.00001074"
Set-TextCell "E18" "  +2.08%  "

Set-TextCell "D19" "0.06599"
Set-TextCell "E19" "  +2.64%  "

Set-TextCell "E20" "  +2.71%  "

Set-TextCell "E21" "  +0.02%  "

Set-TextCell "D22" "6.096"
Set-TextCell "E22" "  +3.88%  "

Set-TextCell "D23" "28.364.23"
Set-TextCell "E23" "  +1.41%  "

Set-TextCell "D24" "11.17"
Set-TextCell "E24" "  +0.36%  "

Set-TextCell "D25" "2.203"
Set-TextCell "E25" "  +5.94%  "

Set-TextCell "D26" "2.456"
Set-TextCell "E26" "  +7.89%  "

Set-TextCell "D27" "157.72"
Set-TextCell "E27" "  +1.25%  "

Set-TextCell "D28" "20.60"
Set-TextCell "E28" "  +2.07%  "

Set-TextCell "D29" "2.038.55"
Set-TextCell "E29" "  +3.11%  "

Set-TextCell "D30" "123.79"
Set-TextCell "E30" "  +3.43%  "

Set-TextCell "E31" "  +1.13%  "

Set-TextCell "E32" "  +4.82%  "

Set-TextCell "D33" "5.650"
Set-TextCell "E33" "  +2.59%  "

Set-TextCell "B34" "Hedera"
Set-TextCell "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D34" "0.07417"
Set-TextCell "E34" "  +16.33%  "

Set-TextCell "B35" "HuobiToken"
Set-TextCell "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D35" "3.648"
Set-TextCell "E35" "  +0.24%  "

Set-TextCell "D36" "0.2232"
Set-TextCell "E36" "  -0.94%  "

Set-TextCell "D37" "0.02344"
Set-TextCell "E37" "  +3.67%  "

Set-TextCell "B38" "FraxShare"
Set-TextCell "C38" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D38" "8.897"
Set-TextCell "E38" "  +6.06%  "

Set-TextCell "B39" "InternetComputer(DFINITY)"
Set-TextCell "C39" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D39" "5.196"
Set-TextCell "E39" "  +4.75%  "

Set-TextCell "B40" "Aptos"
Set-TextCell "C40" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D40" "11.29"
Set-TextCell "E40" "  +2.48%  "

Set-TextCell "B41" "TheSandbox"
Set-TextCell "C41" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D41" "0.6249"
Set-TextCell "E41" "  +2.14%  "

Set-TextCell "D42" "1.183"
Set-TextCell "E42" "  +0.57%  "

Set-TextCell "E43" "  +0.03%  "

Set-TextCell "D44" "1.396"
Set-TextCell "E44" "  -2.25%  "

Set-TextCell "D45" "13.44"
Set-TextCell "E45" "  +1.42%  "

Set-TextCell "D46" "3.697"
Set-TextCell "E46" "  +1.18%  "

Set-TextCell "D47" "0.5839"
Set-TextCell "E47" "  +1.94%  "

Set-TextCell "D48" "124.90"
Set-TextCell "E48" "  -1.07%  "

Set-TextCell "D49" "1.988"
Set-TextCell "E49" "  +3.62%  "

Set-TextCell "D50" "1.201"
Set-TextCell "E50" "  +1.61%  "

Set-TextCell "D51" "0.06893"
Set-TextCell "E51" "  +1.60%  "
